{"js": "// Update the two-digit multiplication equations to their new values.\n// Each old equation string is unique in the document, so a simple\n// search + replace per pair is sufficient and safe.\n\nconst replacements = [\n    { old: \"93\u00d725=2325\", new: \"17\u00d754=918\" },\n    { old: \"34\u00d792=3128\", new: \"91\u00d792=8372\" },\n    { old: \"61\u00d727=1647\", new: \"33\u00d756=1848\" },\n    { old: \"62\u00d785=5270\", new: \"29\u00d768=1972\" },\n    { old: \"57\u00d747=2679\", new: \"16\u00d739=624\" },\n    { old: \"94\u00d764=6016\", new: \"46\u00d737=1702\" },\n    { old: \"33\u00d735=1155\", new: \"88\u00d790=7920\" },\n    { old: \"37\u00d792=3404\", new: \"48\u00d763=3024\" },\n    { old: \"82\u00d794=7708\", new: \"81\u00d761=4941\" },\n    { old: \"47\u00d752=2444\", new: \"78\u00d778=6084\" },\n    { old: \"87\u00d776=6612\", new: \"96\u00d744=4224\" },\n    { old: \"78\u00d793=7254\", new: \"50\u00d799=4950\" },\n    { old: \"49\u00d740=1960\", new: \"20\u00d789=1780\" },\n    { old: \"49\u00d784=4116\", new: \"67\u00d784=5628\" },\n    { old: \"62\u00d752=3224\", new: \"81\u00d755=4455\" },\n    { old: \"81\u00d743=3483\", new: \"39\u00d796=3744\" },\n    { old: \"49\u00d717=833\",  new: \"59\u00d727=1593\" },\n    { old: \"92\u00d727=2484\", new: \"61\u00d717=1037\" },\n    { old: \"19\u00d751=969\",  new: \"22\u00d789=1958\" },\n    { old: \"54\u00d735=1890\", new: \"54\u00d737=1998\" },\n    { old: \"49\u00d743=2107\", new: \"17\u00d772=1224\" },\n    { old: \"40\u00d766=2640\", new: \"91\u00d763=5733\" },\n    { old: \"38\u00d788=3344\", new: \"36\u00d764=2304\" },\n    { old: \"50\u00d736=1800\", new: \"16\u00d740=640\" },\n    { old: \"68\u00d728=1904\", new: \"40\u00d794=3760\" }\n];\n\nconst body = context.document.body;\n\nfor (const pair of replacements) {\n    const results = body.search(pair.old, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(pair.new, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication equations to their new values.\n# Each old equation string is unique in the document, so a simple\n# Find/Replace (wdReplaceAll) per pair is sufficient and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"93\u00d725=2325\"; New = \"17\u00d754=918\" },\n    @{ Old = \"34\u00d792=3128\"; New = \"91\u00d792=8372\" },\n    @{ Old = \"61\u00d727=1647\"; New = \"33\u00d756=1848\" },\n    @{ Old = \"62\u00d785=5270\"; New = \"29\u00d768=1972\" },\n    @{ Old = \"57\u00d747=2679\"; New = \"16\u00d739=624\" },\n    @{ Old = \"94\u00d764=6016\"; New = \"46\u00d737=1702\" },\n    @{ Old = \"33\u00d735=1155\"; New = \"88\u00d790=7920\" },\n    @{ Old = \"37\u00d792=3404\"; New = \"48\u00d763=3024\" },\n    @{ Old = \"82\u00d794=7708\"; New = \"81\u00d761=4941\" },\n    @{ Old = \"47\u00d752=2444\"; New = \"78\u00d778=6084\" },\n    @{ Old = \"87\u00d776=6612\"; New = \"96\u00d744=4224\" },\n    @{ Old = \"78\u00d793=7254\"; New = \"50\u00d799=4950\" },\n    @{ Old = \"49\u00d740=1960\"; New = \"20\u00d789=1780\" },\n    @{ Old = \"49\u00d784=4116\"; New = \"67\u00d784=5628\" },\n    @{ Old = \"62\u00d752=3224\"; New = \"81\u00d755=4455\" },\n    @{ Old = \"81\u00d743=3483\"; New = \"39\u00d796=3744\" },\n    @{ Old = \"49\u00d717=833\";  New = \"59\u00d727=1593\" },\n    @{ Old = \"92\u00d727=2484\"; New = \"61\u00d717=1037\" },\n    @{ Old = \"19\u00d751=969\";  New = \"22\u00d789=1958\" },\n    @{ Old = \"54\u00d735=1890\"; New = \"54\u00d737=1998\" },\n    @{ Old = \"49\u00d743=2107\"; New = \"17\u00d772=1224\" },\n    @{ Old = \"40\u00d766=2640\"; New = \"91\u00d763=5733\" },\n    @{ Old = \"38\u00d788=3344\"; New = \"36\u00d764=2304\" },\n    @{ Old = \"50\u00d736=1800\"; New = \"16\u00d740=640\" },\n    @{ Old = \"68\u00d728=1904\"; New = \"40\u00d794=3760\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
